$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 30,6

$data[0,0] = "A set of methods can be grouped together, ______ shared  variables, into a class.  "
$data[0,1] = "D"
$data[0,2] = "along within some"
$data[0,3] = "along some"
$data[0,4] = "along within"
$data[0,5] = "along with some"

$data[1,0] = "In designing a user interface it is as well to realize  that there are several potentially different viewpoints. The perspectives  include "
$data[1,1] = "D"
$data[1,2] = "the end-user who will eventually get to use the software"
$data[1,3] = "the novice or occasional user"
$data[1,4] = "different end-users with different personalities"
$data[1,5] = "all of the given answers are correct"

$data[2,0] = "In various programming languages, a component is ______ "
$data[2,1] = "A"
$data[2,2] = "all of the given answers are correct"
$data[2,3] = "a method"
$data[2,4] = "a class"
$data[2,5] = "a package"

$data[3,0] = "A novice user or an occasional user______ remember much  about how to use the system.  "
$data[3,1] = "B"
$data[3,2] = "is not"
$data[3,3] = "is not likely to"
$data[3,4] = "are not likely to"
$data[3,5] = "is likely to"

$data[4,0] = "The skill level of the end user has a  significant impact on the ability to  "
$data[4,1] = "A"
$data[4,2] = "all of the given answers are correct"
$data[4,3] = "effectively apply heuristics that create a rhythm of interaction"
$data[4,4] = "respond efficiently to tasks that are demanded by the interaction"
$data[4,5] = "extract meaningful information from the user interface"

$data[5,0] = "Thus a direct manipulation interface______ suitable  approach.  "
$data[5,1] = "B"
$data[5,2] = "be the most"
$data[5,3] = "may be the most"
$data[5,4] = "may be most"
$data[5,5] = "may the most"

$data[6,0] = "The scenario is  software ______ thousands or even hundreds of thousands of lines of code "
$data[6,1] = "C"
$data[6,2] = "that consists off"
$data[6,3] = "that consist of"
$data[6,4] = "that consists of"
$data[6,5] = "that consists in"

$data[7,0] = "It seems that context- or domain-specific knowledge  ______ overall education or intelligence. "
$data[7,1] = ""
$data[7,2] = "is important than"
$data[7,3] = "is more important then"
$data[7,4] = "are more important than"
$data[7,5] = "is more important than"

$data[8,0] = "There is a variety of mechanisms for splitting software  into independent components, or, expressed another way, grouping together items  ______ mutual affinity.  "
$data[8,1] = "D"
$data[8,2] = "that having some"
$data[8,3] = "that has sometimes"
$data[8,4] = "that has some"
$data[8,5] = "that have some"

$data[9,0] = "This aim has consequences _____ stages of software  development, as follows.  "
$data[9,1] = "B"
$data[9,2] = "for nearly that all"
$data[9,3] = "for nearly all"
$data[9,4] = "for all that"
$data[9,5] = "for shall all"

$data[10,0] = "An interface used by two individuals with the  same education and background but entirely different personalities ______ and  unfriendly to the other.  "
$data[10,1] = "A"
$data[10,2] = "may seem friendly to one"
$data[10,3] = "may seem one"
$data[10,4] = "may saw to one"
$data[10,5] = "may see friendly for one"

$data[11,0] = "Therefore, the ideal user interface would be designed to  accommodate differences in personality, or, alternatively, would be designed to  accommodate a typical personality ______.  "
$data[11,1] = "A"
$data[11,2] = "among a class of end users"
$data[11,3] = "a class of end users"
$data[11,4] = "among a class end users"
$data[11,5] = "classes of end users"

$data[12,0] = "At the finest level of granularity, ______  statements and variable declarations can be placed in a method.  "
$data[12,1] = "B"
$data[12,2] = "a numbers of"
$data[12,3] = "a number of"
$data[12,4] = "number of"
$data[12,5] = "a number off"

$data[13,0] = "In designing a user interface it ______ realize that  there are several potentially different viewpoints. "
$data[13,1] = "D"
$data[13,2] = "is as good to"
$data[13,3] = "are good"
$data[13,4] = "is as well for"
$data[13,5] = "is as well to"

$data[14,0] = "The problem is that different people often have  different perspectives of the user interface; they also have different ______.  "
$data[14,1] = "A"
$data[14,2] = "all of the given answers are correct"
$data[14,3] = "skills"
$data[14,4] = "personalities"
$data[14,5] = "culture"

$data[15,0] = "In essence, the desire for modularity ______  construct software from pieces that are as independent of each other as  possible.  "
$data[15,1] = "C"
$data[15,2] = "is trying to be"
$data[15,3] = "is about"
$data[15,4] = "is about trying to"
$data[15,5] = "is to trying to"

$data[16,0] = "For example, a number of applications provide a  macro facility, ______ commands can be grouped together, parameterized and  invoked as a single command "
$data[16,1] = "B"
$data[16,2] = "in which of"
$data[16,3] = "in which a series of"
$data[16,4] = "in series of"
$data[16,5] = "which a series of"

$data[17,0] = "Most people do not apply any formal reasoning ______ a  problem, such as understanding what a computer is displaying.  "
$data[17,1] = "C"
$data[17,2] = "when matched with"
$data[17,3] = "confronting with"
$data[17,4] = "when confronted with"
$data[17,5] = "matching with"

$data[18,0] = "The ______ the desktop metaphor, familiar to users of  Microsoft and Apple Macintosh operating systems.  "
$data[18,1] = "C"
$data[18,2] = "best hidden of those are"
$data[18,3] = "not known of these is"
$data[18,4] = "best known of these is"
$data[18,5] = "best known of these are"

$data[19,0] = "Each person has some model ______ system works  and what it does. "
$data[19,1] = "B"
$data[19,2] = "of whereas the"
$data[19,3] = "of how the"
$data[19,4] = "of who the"
$data[19,5] = "of whom the"

$data[20,0] = "Thus a component is a fairly independent piece of program  ______, some instructions and some data of its own. "
$data[20,1] = "B"
$data[20,2] = "that don’t has a name"
$data[20,3] = "that has a name"
$data[20,4] = "that has name"
$data[20,5] = "that have a name"

$data[21,0] = "These heuristics tend to be domain specific – an  identical problem, encountered in entirely different contexts, ______ applying  different heuristics "
$data[21,1] = "A"
$data[21,2] = "might be solved by"
$data[21,3] = "might be"
$data[21,4] = "might be solving by"
$data[21,5] = "might solve by"

$data[22,0] = "A third possibility ______ interface that is flexible and  can be used in different ways according to personality differences.  "
$data[22,1] = "D"
$data[22,2] = "are to create an"
$data[22,3] = "is create an"
$data[22,4] = "is to create a"
$data[22,5] = "is to create an"

$data[23,0] = "These different perspectives ______ mental models "
$data[23,1] = "C"
$data[23,2] = "is called"
$data[23,3] = "is sometimes called"
$data[23,4] = "are sometimes called"
$data[23,5] = "are sometimes call"

$data[24,0] = "But an experienced and frequent user may be  frustrated by an interface ______ novices and may prefer shortcut commands  and/or a command line interface "
$data[24,1] = "A"
$data[24,2] = "designed for"
$data[24,3] = "designing for"
$data[24,4] = "designed to"
$data[24,5] = "design for"

$data[25,0] = "While there is a massive trend towards multitasking,  window-oriented, point and pick interfaces which can make HCI easier, this  ______ careful design of the interface is conducted. "
$data[25,1] = "D"
$data[25,2] = "only detect"
$data[25,3] = "only happens"
$data[25,4] = "only happen if"
$data[25,5] = "only happens if"

$data[26,0] = "Ideally, each component should be self-contained and  ______ references as possible to other components.  "
$data[26,1] = "D"
$data[26,2] = "has as few"
$data[26,3] = "have few"
$data[26,4] = "has a few"
$data[26,5] = "have as few"

$data[27,0] = "Rather, they apply  a set of ______ based on their understanding of similar problems.  "
$data[27,1] = ""
$data[27,2] = "strategies"
$data[27,3] = "rules"
$data[27,4] = "all of the given answers are correct"
$data[27,5] = "guidelines"

$data[28,0] = "For example, an engineer who uses a computer-based  diagnostic system ______ automobiles understands the problem domain and can  interact effectively through an interface specifically designed to accommodate  users with an engineer’s background. "
$data[28,1] = "C"
$data[28,2] = "to find a fault in"
$data[28,3] = "to finding faults in"
$data[28,4] = "to find faults in"
$data[28,5] = "to find faults"

$data[29,0] = "This same interface might confuse a physician, ______  physician has considerable experience of using a computer for diagnosing  illnesses in patients "
$data[29,1] = "A"
$data[29,2] = "even though the"
$data[29,3] = "the"
$data[29,4] = "even"
$data[29,5] = "though"

$ws.Range("A62:F91").Value = $data

Write-Output "done"